$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Ebook" column header to "BookFormat"
$ws.Range("L1").Value = "BookFormat"

# Update the active selection to K3
$ws.Range("K3").Select()
